# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns for
# rows 2..51 to match the latest scrape.
#
# D-column values are plain text (e.g. "66.235.77", "583.23") mirroring the
# site's formatting (thousand separators as extra dots, fixed 2-decimal
# prices, etc.). Several of them look like ordinary numbers to Excel's
# General-format auto-detection (e.g. "583.23" or "72.90" would silently
# become the numbers 583.23 / 72.9, dropping the trailing zero). To keep
# them as literal text - matching the workbook's existing inline-string
# cells - the D range is switched to Text format ("@") before the values
# are written, then the style is reset back to Normal so no stray
# number-format attribute is left on the cells themselves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value  = "66.235.77"
$ws.Range("E2").Value  = "  -0.38%  "

$ws.Range("D3").Value  = "3.418.17"
$ws.Range("E3").Value  = "  -0.89%  "

$ws.Range("D5").Value  = "583.23"
$ws.Range("E5").Value  = "  -0.29%  "

$ws.Range("D6").Value  = "178.49"
$ws.Range("E6").Value  = "  +1.34%  "

$ws.Range("E7").Value  = "  +3.54%  "

$ws.Range("E8").Value  = "  -0.03%  "

$ws.Range("D9").Value  = "3.414.68"
$ws.Range("E9").Value  = "  -0.90%  "

$ws.Range("E10").Value = "  +0.27%  "

$ws.Range("E11").Value = "  +1.30%  "

$ws.Range("E12").Value = "  -0.77%  "

$ws.Range("D13").Value = "4.010.13"
$ws.Range("E13").Value = "  -0.87%  "

$ws.Range("E14").Value = "  +0.66%  "

$ws.Range("D15").Value = "29.33"
$ws.Range("E15").Value = "  -2.87%  "

$ws.Range("D16").Value = "66.271.77"
$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("E17").Value = "  +0.30%  "

$ws.Range("D18").Value = "3.400.09"
$ws.Range("E18").Value = "  -1.46%  "

$ws.Range("E19").Value = "  -0.66%  "

$ws.Range("D20").Value = "13.76"
$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("D21").Value = "366.61"
$ws.Range("E21").Value = "  -3.11%  "

$ws.Range("D22").Value = "7.58"
$ws.Range("E22").Value = "  -3.07%  "

$ws.Range("D23").Value = "72.90"
$ws.Range("E23").Value = "  +0.76%  "

$ws.Range("E24").Value = "  -0.35%  "

$ws.Range("E25").Value = "  +5.27%  "

$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("E28").Value = "  +1.99%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").Value = "5.77"
$ws.Range("E30").Value = "  -1.21%  "

$ws.Range("E31").Value = "  -0.56%  "

$ws.Range("D32").Value = "23.41"
$ws.Range("E32").Value = "  -3.54%  "

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").Value = "7.01"
$ws.Range("E34").Value = "  -2.13%  "

$ws.Range("E35").Value = "  -3.92%  "

$ws.Range("E36").Value = "  -1.53%  "

$ws.Range("D37").Value = "162.94"
$ws.Range("E37").Value = "  +1.11%  "

$ws.Range("D38").Value = "0.873"
$ws.Range("E38").Value = "  -1.94%  "

$ws.Range("D39").Value = "27.60"
$ws.Range("E39").Value = "  -5.93%  "

$ws.Range("E40").Value = "  +0.37%  "

$ws.Range("E41").Value = "  -1.87%  "

$ws.Range("E42").Value = "  -1.19%  "

$ws.Range("D43").Value = "2.698.99"
$ws.Range("E43").Value = "  -1.41%  "

$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("E45").Value = "  -0.98%  "

$ws.Range("D46").Value = "25.20"
$ws.Range("E46").Value = "  +2.83%  "

$ws.Range("D47").Value = "39.92"
$ws.Range("E47").Value = "  -1.81%  "

$ws.Range("D48").Value = "336.41"
$ws.Range("E48").Value = "  +8.95%  "

$ws.Range("E49").Value = "  -2.74%  "

$ws.Range("E50").Value = "  +2.37%  "

$ws.Range("D51").Value = "31.99"
$ws.Range("E51").Value = "  +4.37%  "

# Restore the default style on the D range so cells don't carry a stray
# explicit style index - only their stored value/type changed.
$dRange.Style = "Normal"
